$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REGISTRU")
$row = 7

# Copy the formatting (border/alignment/bold "ID" style) from the row above
# into the new ID cell before writing its value.
$ws.Range("A6").Copy()
$ws.Range("A" + $row).PasteSpecial(-4122)

# CNP / TELEFON / NUMAR_BILET / PRESIUNE look numeric but must stay text,
# exactly like every other data row in the sheet.
$ws.Cells.Item($row, 2).NumberFormat  = "@"
$ws.Cells.Item($row, 5).NumberFormat  = "@"
$ws.Cells.Item($row, 6).NumberFormat  = "@"
$ws.Cells.Item($row, 13).NumberFormat = "@"
$ws.Cells.Item($row, 19).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "03-01-2024"
$ws.Cells.Item($row, 3).Value  = "LEO"
$ws.Cells.Item($row, 4).Value  = "MESSI"
$ws.Cells.Item($row, 5).Value  = "1900106226823"
$ws.Cells.Item($row, 6).Value  = "5522333"
$ws.Cells.Item($row, 7).Value  = "INTER"
$ws.Cells.Item($row, 8).Value  = "MIMAI"
$ws.Cells.Item($row, 9).Value  = "STRAINATATE"
$ws.Cells.Item($row, 10).Value = "YES"
$ws.Cells.Item($row, 11).Value = "Salariat"
$ws.Cells.Item($row, 12).Value = "YES"
$ws.Cells.Item($row, 13).Value = "256314"
$ws.Cells.Item($row, 14).Value = "JUCATOR BUN`n"
$ws.Cells.Item($row, 15).Value = "YES"
$ws.Cells.Item($row, 16).Value = "OBSTRUCTIV"
$ws.Cells.Item($row, 17).Value = "PILLOWS"
$ws.Cells.Item($row, 18).Value = "BUNA"
$ws.Cells.Item($row, 19).Value = "0.325"
$ws.Cells.Item($row, 20).Value = "YES"
$ws.Cells.Item($row, 21).Value = "TALENT PREA MARE`n"
$ws.Cells.Item($row, 22).Value = "A JUCA PANA LA CUPA MONDIALA SI A O CASTIGA`n"

# Grow the AutoFilter range to include the new row (A1:V6 -> A1:V7).
$ws.AutoFilterMode = $false
$ws.Range("A1:V7").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "REGISTRU!_FilterDatabase") {
        $n.RefersTo = "='REGISTRU'!`$A`$1:`$V`$7"
    }
}
